$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10
$ws.Range("I10").Value = "sd"
$ws.Range("J10").Value = "Statement-non-opinion"

# Row 15
$ws.Range("I15").Value = "b"
$ws.Range("J15").Value = "Acknowledge (Backchannel)"

# Row 22
$ws.Range("I22").Value = "b"
$ws.Range("J22").Value = "Acknowledge (Backchannel)"

# Row 31
$ws.Range("I31").Value = "sd"
$ws.Range("J31").Value = "Statement-non-opinion"

# Row 40
$ws.Range("I40").Value = "aa"
$ws.Range("J40").Value = "Agree/Accept"

# Row 48
$ws.Range("I48").Value = "sv"
$ws.Range("J48").Value = "Statement-opinion"

# Row 54
$ws.Range("I54").Value = "sd"
$ws.Range("J54").Value = "Statement-non-opinion"

# Row 69
$ws.Range("I69").Value = "sv"
$ws.Range("J69").Value = "Statement-opinion"

# Row 72
$ws.Range("I72").Value = "b"
$ws.Range("J72").Value = "Acknowledge (Backchannel)"

# Row 84
$ws.Range("I84").Value = "sd"
$ws.Range("J84").Value = "Statement-non-opinion"

# Row 87
$ws.Range("I87").Value = "sv"
$ws.Range("J87").Value = "Statement-opinion"
